$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# 1) Column D (bsecode) for rows 158-172: convert from text to numeric, same values
$ws.Cells.Item(158, 4).Value = 532538
$ws.Cells.Item(159, 4).Value = 539523
$ws.Cells.Item(160, 4).Value = 505200
$ws.Cells.Item(161, 4).Value = 500410
$ws.Cells.Item(162, 4).Value = 532830
$ws.Cells.Item(163, 4).Value = 500483
$ws.Cells.Item(164, 4).Value = 539957
$ws.Cells.Item(165, 4).Value = 500575
$ws.Cells.Item(166, 4).Value = 543066
$ws.Cells.Item(167, 4).Value = 500425
$ws.Cells.Item(168, 4).Value = 531642
$ws.Cells.Item(169, 4).Value = 539336
$ws.Cells.Item(170, 4).Value = 512070
$ws.Cells.Item(171, 4).Value = 533519
$ws.Cells.Item(172, 4).Value = 540750

# 2) Append 15 new rows (173-187) duplicating the 158-172 block with updated stats
# Column D stays text (bsecode as string) in these new rows, per source diff
$ws.Cells.Item(173, 1).Value = 1
$ws.Cells.Item(173, 2).Value = "ULTRACEMCO"
$ws.Cells.Item(173, 3).Value = "Ultratech Cement Limited"
$ws.Cells.Item(173, 4).NumberFormat = "@"
$ws.Cells.Item(173, 4).Value = "532538"
$ws.Cells.Item(173, 4).Style = "Normal"
$ws.Cells.Item(173, 5).Value = -1.48
$ws.Cells.Item(173, 6).Value = 11658.7
$ws.Cells.Item(173, 7).Value = 556056
$ws.Cells.Item(173, 8).Value = "day"
$ws.Cells.Item(173, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(174, 1).Value = 2
$ws.Cells.Item(174, 2).Value = "ALKEM"
$ws.Cells.Item(174, 3).Value = "Alkem Laboratories Limited"
$ws.Cells.Item(174, 4).NumberFormat = "@"
$ws.Cells.Item(174, 4).Value = "539523"
$ws.Cells.Item(174, 4).Style = "Normal"
$ws.Cells.Item(174, 5).Value = -1.61
$ws.Cells.Item(174, 6).Value = 5296.8
$ws.Cells.Item(174, 7).Value = 229725
$ws.Cells.Item(174, 8).Value = "day"
$ws.Cells.Item(174, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "EICHERMOT"
$ws.Cells.Item(175, 3).Value = "Eicher Motors Limited"
$ws.Cells.Item(175, 4).NumberFormat = "@"
$ws.Cells.Item(175, 4).Value = "505200"
$ws.Cells.Item(175, 4).Style = "Normal"
$ws.Cells.Item(175, 5).Value = 0.71
$ws.Cells.Item(175, 6).Value = 4916.1
$ws.Cells.Item(175, 7).Value = 664540
$ws.Cells.Item(175, 8).Value = "day"
$ws.Cells.Item(175, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(176, 1).Value = 4
$ws.Cells.Item(176, 2).Value = "ACC"
$ws.Cells.Item(176, 3).Value = "Acc Limited"
$ws.Cells.Item(176, 4).NumberFormat = "@"
$ws.Cells.Item(176, 4).Value = "500410"
$ws.Cells.Item(176, 4).Style = "Normal"
$ws.Cells.Item(176, 5).Value = 0.75
$ws.Cells.Item(176, 6).Value = 2715.85
$ws.Cells.Item(176, 7).Value = 446114
$ws.Cells.Item(176, 8).Value = "day"
$ws.Cells.Item(176, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(177, 1).Value = 5
$ws.Cells.Item(177, 2).Value = "ASTRAL"
$ws.Cells.Item(177, 3).Value = "Astral Poly Technik Limited"
$ws.Cells.Item(177, 4).NumberFormat = "@"
$ws.Cells.Item(177, 4).Value = "532830"
$ws.Cells.Item(177, 4).Style = "Normal"
$ws.Cells.Item(177, 5).Value = 0.19
$ws.Cells.Item(177, 6).Value = 2266.85
$ws.Cells.Item(177, 7).Value = 210636
$ws.Cells.Item(177, 8).Value = "day"
$ws.Cells.Item(177, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(178, 1).Value = 6
$ws.Cells.Item(178, 2).Value = "TATACOMM"
$ws.Cells.Item(178, 3).Value = "Tata Communications Limited"
$ws.Cells.Item(178, 4).NumberFormat = "@"
$ws.Cells.Item(178, 4).Value = "500483"
$ws.Cells.Item(178, 4).Style = "Normal"
$ws.Cells.Item(178, 5).Value = -1.07
$ws.Cells.Item(178, 6).Value = 1851.05
$ws.Cells.Item(178, 7).Value = 223941
$ws.Cells.Item(178, 8).Value = "day"
$ws.Cells.Item(178, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(179, 1).Value = 7
$ws.Cells.Item(179, 2).Value = "MGL"
$ws.Cells.Item(179, 3).Value = "Mahanagar Gas Limited"
$ws.Cells.Item(179, 4).NumberFormat = "@"
$ws.Cells.Item(179, 4).Value = "539957"
$ws.Cells.Item(179, 4).Style = "Normal"
$ws.Cells.Item(179, 5).Value = -0.3
$ws.Cells.Item(179, 6).Value = 1736.95
$ws.Cells.Item(179, 7).Value = 337792
$ws.Cells.Item(179, 8).Value = "day"
$ws.Cells.Item(179, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(180, 1).Value = 8
$ws.Cells.Item(180, 2).Value = "VOLTAS"
$ws.Cells.Item(180, 3).Value = "Voltas Limited"
$ws.Cells.Item(180, 4).NumberFormat = "@"
$ws.Cells.Item(180, 4).Value = "500575"
$ws.Cells.Item(180, 4).Style = "Normal"
$ws.Cells.Item(180, 5).Value = -0.6
$ws.Cells.Item(180, 6).Value = 1521.7
$ws.Cells.Item(180, 7).Value = 690677
$ws.Cells.Item(180, 8).Value = "day"
$ws.Cells.Item(180, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(181, 1).Value = 9
$ws.Cells.Item(181, 2).Value = "SBICARD"
$ws.Cells.Item(181, 3).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(181, 4).NumberFormat = "@"
$ws.Cells.Item(181, 4).Value = "543066"
$ws.Cells.Item(181, 4).Style = "Normal"
$ws.Cells.Item(181, 5).Value = -1.04
$ws.Cells.Item(181, 6).Value = 730.9
$ws.Cells.Item(181, 7).Value = 1032022
$ws.Cells.Item(181, 8).Value = "day"
$ws.Cells.Item(181, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(182, 1).Value = 10
$ws.Cells.Item(182, 2).Value = "AMBUJACEM"
$ws.Cells.Item(182, 3).Value = "Ambuja Cements Limited"
$ws.Cells.Item(182, 4).NumberFormat = "@"
$ws.Cells.Item(182, 4).Value = "500425"
$ws.Cells.Item(182, 4).Style = "Normal"
$ws.Cells.Item(182, 5).Value = 0.17
$ws.Cells.Item(182, 6).Value = 685.35
$ws.Cells.Item(182, 7).Value = 2484802
$ws.Cells.Item(182, 8).Value = "day"
$ws.Cells.Item(182, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(183, 1).Value = 11
$ws.Cells.Item(183, 2).Value = "MARICO"
$ws.Cells.Item(183, 3).Value = "Marico Limited"
$ws.Cells.Item(183, 4).NumberFormat = "@"
$ws.Cells.Item(183, 4).Value = "531642"
$ws.Cells.Item(183, 4).Style = "Normal"
$ws.Cells.Item(183, 5).Value = 2.21
$ws.Cells.Item(183, 6).Value = 667.35
$ws.Cells.Item(183, 7).Value = 3414973
$ws.Cells.Item(183, 8).Value = "day"
$ws.Cells.Item(183, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(184, 1).Value = 12
$ws.Cells.Item(184, 2).Value = "GUJGASLTD"
$ws.Cells.Item(184, 3).Value = "Gujarat Gas Limited"
$ws.Cells.Item(184, 4).NumberFormat = "@"
$ws.Cells.Item(184, 4).Value = "539336"
$ws.Cells.Item(184, 4).Style = "Normal"
$ws.Cells.Item(184, 5).Value = -0.73
$ws.Cells.Item(184, 6).Value = 637.8
$ws.Cells.Item(184, 7).Value = 1356224
$ws.Cells.Item(184, 8).Value = "day"
$ws.Cells.Item(184, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(185, 1).Value = 13
$ws.Cells.Item(185, 2).Value = "UPL"
$ws.Cells.Item(185, 3).Value = "Upl Limited"
$ws.Cells.Item(185, 4).NumberFormat = "@"
$ws.Cells.Item(185, 4).Value = "512070"
$ws.Cells.Item(185, 4).Style = "Normal"
$ws.Cells.Item(185, 5).Value = -1.37
$ws.Cells.Item(185, 6).Value = 557.3
$ws.Cells.Item(185, 7).Value = 1674706
$ws.Cells.Item(185, 8).Value = "day"
$ws.Cells.Item(185, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(186, 1).Value = 14
$ws.Cells.Item(186, 2).Value = "LTF"
$ws.Cells.Item(186, 3).Value = "L&T Finance Ltd"
$ws.Cells.Item(186, 4).NumberFormat = "@"
$ws.Cells.Item(186, 4).Value = "533519"
$ws.Cells.Item(186, 4).Style = "Normal"
$ws.Cells.Item(186, 5).Value = -0.81
$ws.Cells.Item(186, 6).Value = 184.35
$ws.Cells.Item(186, 7).Value = 12008676
$ws.Cells.Item(186, 8).Value = "day"
$ws.Cells.Item(186, 9).Value = "17/07/2024 11:34:50"

$ws.Cells.Item(187, 1).Value = 15
$ws.Cells.Item(187, 2).Value = "IEX"
$ws.Cells.Item(187, 3).Value = "Indian Energy Exchange Ltd"
$ws.Cells.Item(187, 4).NumberFormat = "@"
$ws.Cells.Item(187, 4).Value = "540750"
$ws.Cells.Item(187, 4).Style = "Normal"
$ws.Cells.Item(187, 5).Value = -0.49
$ws.Cells.Item(187, 6).Value = 177.34
$ws.Cells.Item(187, 7).Value = 11969269
$ws.Cells.Item(187, 8).Value = "day"
$ws.Cells.Item(187, 9).Value = "17/07/2024 11:34:50"

